$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 240
$ws.Range("D3").Value = 240

$ws.Range("D2:D3").Select()
